# Auto-generated Excel COM-interop script applying the Hyperion_Profits.xlsx diff.
# Each row below mirrors a single <row> edit from the unified OOXML diff:
#   - "change": existing cell value updated to a new number
#   - "added":  a previously-empty cell now holds a value
#   - "removed": a previously-populated cell is cleared (no <c> element remains)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100
$ws.Range("H100").Value = 600
$ws.Range("I100").Value = 600
$ws.Range("K100").Value = 600
$ws.Range("M100").Value = -59

# Row 125
$ws.Range("H125").Value = 9011550
$ws.Range("I125").Value = 1586.8334
$ws.Range("J125").Value = 10755413
$ws.Range("K125").Value = 14281.5006
$ws.Range("L125").Value = 96798717
$ws.Range("M125").Value = -11821.5006
$ws.Range("N125").Value = -96803637

# Row 132
$ws.Range("H132").Value = 55559924
$ws.Range("I132").Value = 62504540
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 187513620
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -187511090
$ws.Range("N132").Value = -14058.5

# Row 137
$ws.Range("H137").Value = 120818.2
$ws.Range("I137").Value = 179398
$ws.Range("J137").Value = 3658.6
$ws.Range("K137").Value = 538194
$ws.Range("L137").Value = 10975.8
$ws.Range("M137").Value = -535644
$ws.Range("N137").Value = -16075.8

# Row 138
$ws.Range("H138").Value = 2180.5422
$ws.Range("J138").Value = 3480.4
$ws.Range("L138").Value = 10441.2
$ws.Range("N138").Value = -20721.2

# Row 141
$ws.Range("H141").Value = 1932.909
$ws.Range("I141").Value = 1729.625
$ws.Range("J141").Value = 2475
$ws.Range("K141").Value = 5188.875
$ws.Range("L141").Value = 7425
$ws.Range("M141").Value = -8.875
$ws.Range("N141").Value = -17785

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 1208.3334
$ws.Range("I5").Value = 1045.4546
$ws.Range("K5").Value = 1045.4546
$ws.Range("M5").Value = -933.4546

# Row 32
$ws.Range("H32").Value = 5327.981
$ws.Range("I32").Value = 2993.425
$ws.Range("K32").Value = 2993.425
$ws.Range("M32").Value = -2706.425

# Row 45
$ws.Range("H45").Value = 6996761
$ws.Range("I45").Value = 11836717
$ws.Range("J45").Value = 5713
$ws.Range("K45").Value = 11836717
$ws.Range("L45").Value = 5713
$ws.Range("M45").Value = -11836340
$ws.Range("N45").Value = -6467

# Row 102
$ws.Range("H102").Value = 4633220.5
$ws.Range("I102").Value = 6946929
$ws.Range("K102").Value = 6946929
$ws.Range("M102").Value = -6945307

# Row 122
$ws.Range("I122").Value = 2136.3928
$ws.Range("J122").Value = 2086472.9
$ws.Range("K122").Value = 6409.178400000001
$ws.Range("L122").Value = 6259418.699999999
$ws.Range("M122").Value = -3959.178400000001
$ws.Range("N122").Value = -6264318.699999999

# Row 132
$ws.Range("H132").Value = 3102.3
$ws.Range("I132").Value = 2607.182
$ws.Range("J132").Value = 3388.9473
$ws.Range("K132").Value = 7821.545999999999
$ws.Range("L132").Value = 10166.8419
$ws.Range("M132").Value = -5291.545999999999
$ws.Range("N132").Value = -15226.8419

# Row 133
$ws.Range("H133").Value = 40000
$ws.Range("I133").Value = 40000
$ws.Range("K133").Value = 40000
$ws.Range("M133").Value = -37470

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 1208.3334
$ws.Range("I4").Value = 1045.4546
$ws.Range("K4").Value = 1045.4546
$ws.Range("M4").Value = -930.4546

# Row 107
$ws.Range("H107").Value = 7149456.5
$ws.Range("I107").Value = 14289914
$ws.Range("J107").Value = 8999.6
$ws.Range("K107").Value = 14289914
$ws.Range("L107").Value = 8999.6
$ws.Range("M107").Value = -14287994
$ws.Range("N107").Value = -12839.6

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2177
$ws.Range("I16").Value = 1722.5
$ws.Range("J16").Value = 3995
$ws.Range("K16").Value = 1722.5
$ws.Range("L16").Value = 3995
$ws.Range("M16").Value = -1435.5
$ws.Range("N16").Value = -4569

# Row 31
$ws.Range("H31").Value = 14603.606
$ws.Range("I31").Value = 2443.6
$ws.Range("J31").Value = 16142.848
$ws.Range("K31").Value = 2443.6
$ws.Range("L31").Value = 16142.848
$ws.Range("M31").Value = -2148.6
$ws.Range("N31").Value = -16732.848

# Row 34
$ws.Range("H34").Value = 14603.606
$ws.Range("I34").Value = 2443.6
$ws.Range("J34").Value = 16142.848
$ws.Range("K34").Value = 2443.6
$ws.Range("L34").Value = 16142.848
$ws.Range("M34").Value = -2241.6
$ws.Range("N34").Value = -16546.848

# Row 58
$ws.Range("H58").Value = 3335.8333
$ws.Range("I58").Value = 2570.5715
$ws.Range("J58").Value = 4407.2
$ws.Range("K58").Value = 2570.5715
$ws.Range("L58").Value = 4407.2
$ws.Range("M58").Value = -2367.5715
$ws.Range("N58").Value = -4813.2

# Row 107
$ws.Range("H107").Value = 1407.138
$ws.Range("I107").Value = 1379.9546
$ws.Range("K107").Value = 1379.9546
$ws.Range("M107").Value = 540.0454

# Row 113
$ws.Range("H113").Value = 2177
$ws.Range("I113").Value = 1722.5
$ws.Range("J113").Value = 3995
$ws.Range("K113").Value = 1722.5
$ws.Range("L113").Value = 3995
$ws.Range("M113").Value = 447.5
$ws.Range("N113").Value = -8335

# Row 122
$ws.Range("H122").Value = 4831.8335
$ws.Range("I122").Value = 4165.3335
$ws.Range("K122").Value = 12496.0005
$ws.Range("M122").Value = -10046.0005

# Row 127
$ws.Range("H127").Value = 83299.664
$ws.Range("J127").Value = 83299.664
$ws.Range("L127").Value = 83299.664
$ws.Range("N127").Value = -93219.664

# Row 132
$ws.Range("H132").Value = 49103.895
$ws.Range("I132").Value = 2449.2222
$ws.Range("K132").Value = 7347.6666
$ws.Range("M132").Value = -4817.6666

# Row 134
$ws.Range("H134").Value = 2754.16
$ws.Range("I134").Value = 1892.6471
$ws.Range("K134").Value = 5677.9413
$ws.Range("M134").Value = -3142.9413

# Row 136
$ws.Range("H136").Value = 3335.8333
$ws.Range("I136").Value = 2570.5715
$ws.Range("J136").Value = 4407.2
$ws.Range("K136").Value = 7711.7145
$ws.Range("L136").Value = 13221.6
$ws.Range("M136").Value = -5161.7145
$ws.Range("N136").Value = -18321.6

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 4377588
$ws.Range("I4").Value = 6704147.5
$ws.Range("K4").Value = 20112442.5
$ws.Range("M4").Value = -20112330.5

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 7235.857
$ws.Range("I2").Value = 117.2
$ws.Range("J2").Value = 25032.5
$ws.Range("K2").Value = 117.2
$ws.Range("L2").Value = 25032.5
$ws.Range("M2").Value = -4.200000000000003
$ws.Range("N2").Value = -25258.5

# Row 20
$ws.Range("H20").Value = 5003602.5
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

# Row 25
$ws.Range("H25").Value = 1000
$ws.Range("J25").Value = 1000
$ws.Range("L25").Value = 1000
$ws.Range("N25").Value = -2058

# Row 102
$ws.Range("H102").Value = 12269000
$ws.Range("I102").Value = 22224698
$ws.Range("J102").Value = 3972585.8
$ws.Range("K102").Value = 22224698
$ws.Range("L102").Value = 3972585.8
$ws.Range("M102").Value = -22223076
$ws.Range("N102").Value = -3975829.8

# Row 122
$ws.Range("H122").Value = 273076.34
$ws.Range("I122").Value = 389503.56
$ws.Range("K122").Value = 1168510.68
$ws.Range("M122").Value = -1166060.68

# Row 124
$ws.Range("H124").Value = 42000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 42000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 42000
$ws.Range("M124").Value = ""
$ws.Range("N124").Value = -51820

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4995.161
$ws.Range("I40").Value = 4456.2856
$ws.Range("J40").Value = 5438.9414
$ws.Range("K40").Value = 4456.2856
$ws.Range("L40").Value = 5438.9414
$ws.Range("M40").Value = -4320.2856
$ws.Range("N40").Value = -5710.9414

# Row 46
$ws.Range("H46").Value = 4498.9287
$ws.Range("I46").Value = 926.8570999999999
$ws.Range("J46").Value = 8071
$ws.Range("K46").Value = 926.8570999999999
$ws.Range("L46").Value = 8071
$ws.Range("M46").Value = -738.8570999999999
$ws.Range("N46").Value = -8447

# Row 61
$ws.Range("H61").Value = 6175944
$ws.Range("I61").Value = 6947825
$ws.Range("K61").Value = 6947825
$ws.Range("M61").Value = -6947623

# Row 113
$ws.Range("H113").Value = 6175944
$ws.Range("I113").Value = 6947825
$ws.Range("K113").Value = 6947825
$ws.Range("M113").Value = -6945655

# Row 122
$ws.Range("H122").Value = 5246.7915
$ws.Range("I122").Value = 3831.7646
$ws.Range("K122").Value = 11495.2938
$ws.Range("M122").Value = -9045.293799999999

# Row 136
$ws.Range("H136").Value = 28102.924
$ws.Range("I136").Value = 37732.93
$ws.Range("J136").Value = 3590.182
$ws.Range("K136").Value = 113198.79
$ws.Range("L136").Value = 10770.546
$ws.Range("M136").Value = -110648.79
$ws.Range("N136").Value = -15870.546

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 49000
$ws.Range("J64").Value = 49000
$ws.Range("L64").Value = 49000
$ws.Range("N64").Value = -49496

# Row 67
$ws.Range("H67").Value = 49000
$ws.Range("J67").Value = 49000
$ws.Range("L67").Value = 49000
$ws.Range("N67").Value = -50716

# Row 107
$ws.Range("H107").Value = 55558324
$ws.Range("I107").Value = 62501180
$ws.Range("K107").Value = 187503540
$ws.Range("M107").Value = -187501620

# Row 122
$ws.Range("H122").Value = 2157.8572
$ws.Range("I122").Value = 1422
$ws.Range("J122").Value = 3997.5
$ws.Range("K122").Value = 4266
$ws.Range("L122").Value = 11992.5
$ws.Range("M122").Value = -1816
$ws.Range("N122").Value = -16892.5

# Row 126
$ws.Range("H126").Value = 4244.8184
$ws.Range("I126").Value = 3961.625
$ws.Range("K126").Value = 11884.875
$ws.Range("M126").Value = -9414.875

# Row 132
$ws.Range("H132").Value = 34519330
$ws.Range("I132").Value = 58831980
$ws.Range("J132").Value = 76398.414
$ws.Range("K132").Value = 176495940
$ws.Range("L132").Value = 229195.242
$ws.Range("M132").Value = -176493410
$ws.Range("N132").Value = -234255.242
